# UndoRedoSequenceDiagram.pptx update
# - "popUndo()"          -> "undo()"
# - "x:XYZCommand"       -> ":" + "ModelManager"
# - ":UndoRedoStack"     -> ":" + "UndoRedoCareTaker"  (widen the box)
# - "undo()"             -> "resetData" + "(" + "AddressBook" + ")"  (move/resize box)
#
# NB: Shape.Left/.Top/.Width/.Height are expressed in points (1 pt = 12700 EMU),
# matching the real PowerPoint object model.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById {
    param($slide, [int]$id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) { return $candidate }
    }
    return $null
}

# 1) "popUndo()" -> "undo()"   (shape id 79, "TextBox 78")
#    Original runs: "popUndo" (err=1) + "()".  Empty out the first run so the
#    surviving (clean) run absorbs the whole replacement text -- this avoids
#    PowerPoint's "preserve unchanged text formatting" diff merging the new
#    text across both of the old runs.
$shpPopUndo = Get-ShapeById $s 79
$trPopUndo = $shpPopUndo.TextFrame.TextRange
$trPopUndo.Characters(1, 7).Text = ""
$trPopUndo.Text = "undo()"

# 2) "x:XYZCommand" -> ":" + "ModelManager"   (shape id 84, "Rectangle 62")
$shpXyz = Get-ShapeById $s 84
$trXyz = $shpXyz.TextFrame.TextRange
$lenXyz = $trXyz.Length
$trXyz.Characters(2, $lenXyz - 1).Text = "ModelManager"
$trXyz.Characters(1, 1).Text = ":"

# 3) ":UndoRedoStack" -> ":" + "UndoRedoCareTaker", widen the box  (shape id 40)
$shpCareTaker = Get-ShapeById $s 40
$shpCareTaker.Width = 156.0748071496063   # 1982150 EMU
$trCareTaker = $shpCareTaker.TextFrame.TextRange
$lenCareTaker = $trCareTaker.Length
$trCareTaker.Characters(2, $lenCareTaker - 1).Text = "UndoRedoCareTaker"

# 4) "undo()" -> "resetData" + "(" + "AddressBook" + ")", reposition/widen box
#    (shape id 88, "TextBox 87")
$shpReset = Get-ShapeById $s 88
$shpReset.Left = 483.6592255984252     # 6142472 EMU
$shpReset.Top = 269.0087441574803      # 3416411 EMU
$shpReset.Width = 127.19905511811024   # 1615428 EMU
$trReset = $shpReset.TextFrame.TextRange
$trReset.Text = "resetData(AddressBook)"
$trReset.Characters(22, 1).Text = ")"
$trReset.Characters(11, 11).Text = "AddressBook"
$trReset.Characters(10, 1).Text = "("
$trReset.Characters(1, 9).Text = "resetData"
